# Refresh cryptocurrency price / 1h-volume data (GitHub Actions job).
# Price/volume cells hold plain text in the workbook. Values that look
# like plain numbers (e.g. "1.00", "26.30") are written via a quick
# Text-format round-trip so Excel keeps them as text instead of
# re-interpreting them as numbers; the number-format tweak is cleared
# right after so the cell keeps its original (default/General) style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.082.94"
$ws.Range("E2").Value = "  -2.98%  "
$ws.Range("D3").Value = "3.839.60"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.81"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.50"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.38%  "
$ws.Range("D7").Value = "3.832.98"
$ws.Range("E7").Value = "  -2.56%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -1.89%  "
$ws.Range("E10").Value = "  -4.10%  "
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.71%  "
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.12"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.00%  "
$ws.Range("D15").Value = "4.486.20"
$ws.Range("E15").Value = "  -2.44%  "
$ws.Range("D16").Value = "3.839.65"
$ws.Range("E16").Value = "  -2.84%  "
$ws.Range("D17").Value = "68.177.38"
$ws.Range("E17").Value = "  -2.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.35"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.40%  "
$ws.Range("E19").Value = "  -3.76%  "
$ws.Range("E20").Value = "  -0.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.12"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "466.85"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -6.14%  "
$ws.Range("E23").Value = "  -1.95%  "
$ws.Range("E24").Value = "  -3.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.93"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.68%  "
$ws.Range("E26").Value = "  -3.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.13"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.06"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.67%  "
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("E30").Value = "  -1.69%  "
$ws.Range("D31").Value = "3.990.51"
$ws.Range("E31").Value = "  -2.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.66"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.14%  "
$ws.Range("E33").Value = "  -6.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.30"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.58"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.84%  "
$ws.Range("D36").Value = "3.803.25"
$ws.Range("E36").Value = "  -2.51%  "
$ws.Range("E37").Value = "  -3.58%  "
$ws.Range("E38").Value = "  +10.30%  "
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("E40").Value = "  -2.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.93"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.315"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.87%  "
$ws.Range("E44").Value = "  -7.30%  "
$ws.Range("E45").Value = "  -3.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.71"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000292"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +5.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "46.92"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.99%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "26.30"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.96%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "142.11"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.94%  "
